$wb = $excel.ActiveWorkbook

# --- "Fakers" sheet: update faker-check statuses and add newly processed rows ---
$ws = $wb.Worksheets.Item("Fakers")

# Rows 50-59: "Ответ" (answer) column flips from "Нет" to "Да"
$ws.Range("C50:C59").Value = "Да"
# Rows 57-59: "Блок" (block) column flips from "-" to "Да"
$ws.Range("D57:D59").Value = "Да"

# New row 60: Disha
$ws.Range("A60").Value = "Disha"
$ws.Range("B60").Value = 17813456
$ws.Range("C60").Value = "Да"
$ws.Range("D60").Value = "Да"
$ws.Hyperlinks.Add($ws.Range("E60"), "https://www.ozon.ru/seller/disha-167928/", "", "", "https://www.ozon.ru/seller/disha-167928/")
$ws.Range("E59").Copy()
$ws.Range("E60").PasteSpecial(-4122)

# New row 61: УЮТНЫЙ ДОМ
$ws.Range("A61").Value = "УЮТНЫЙ ДОМ"
$ws.Range("B61").Value = 17813535
$ws.Range("C61").Value = "Да"
$ws.Range("D61").Value = "Да"
$ws.Range("E61").Value = "https://www.ozon.ru/seller/uyutnyy-dom-271916/"

# New row 62: Бейсболка
$ws.Range("A62").Value = "Бейсболка"
$ws.Range("B62").Value = 17813484
$ws.Range("C62").Value = "Да"
$ws.Range("D62").Value = "Да"
$ws.Range("E62").Value = "https://www.ozon.ru/seller/beysbolka-466162/"

# --- "Проверка" sheet: re-enter the VLOOKUP formulas as one pass so Excel groups them as a shared formula ---
$ws3 = $wb.Worksheets.Item("Проверка")
$ws3.Range("B2:B15").Formula = "=VLOOKUP(A2,C:C,1,0)"

# --- Switch the active/selected sheet & view back to "Fakers" ---
$ws.Activate() | Out-Null
$ws.Range("D63").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
